# FHIR-26840 remove prov reqs
# The "Medications" row previously listed three US Core profiles/resources
# (Medication Profile, Medication Statement Profile, Medication Request
# Profile / Medication, MedicationStatement, MedicationRequest). The edit
# drops the "Medication Statement" profile+resource from both the summary
# table (Sheet1) and the crosswalk table (Sheet2).

$wb = $excel.ActiveWorkbook

# --- Sheet1: "USCDI v1 Summary of Data Classes and Data Elements" table ---
# Row 23 is the "  Medications" line -> update its Resource (C) column
# first, then the crosswalk sheet, then the Profile (B) column -- matching
# the order the new shared strings were authored in.
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

$ws1.Range("C23").Value = "Medication,  MedicationRequest"

# --- Sheet2: crosswalk table ---
# Row 11 holds the combined tab-separated "Medications" row.
$ws2.Range("A11").Value = "Medications`t[US Core Medication Profile], [US Core Medication Reqsuest Profile]`tMedication, MedicationRequest"

$ws1.Range("B23").Value = "[US Core Medication Profile], [US Core Medication Request Profile]"

# --- Restore view/selection state ---
# Touch Sheet2's selection first, then return focus to Sheet1 so it stays
# the active tab (matches the saved sheetViews in the workbook).
$ws2.Activate()
$ws2.Range("A11").Select()

$ws1.Activate()
$ws1.Range("B14").Select()
